$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.720118880271912
$ws.Range("B1").Value = 2.714961528778076
$ws.Range("C1").Value = 2.921954154968262
$ws.Range("D1").Value = 3.294327020645142
$ws.Range("E1").Value = 2.475978136062622
